$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.359.20'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.30%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.686.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.42%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.012'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '218.38'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.41%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5441'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.64%  '
$ws.Range("E7").Value = '  +0.72%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2724'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.05%  '
$ws.Range("E9").Value = '  -0.11%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.98'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.97%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07683'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.75%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.700.79'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.530'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.5809'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000008370'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -2.60%  '
$ws.Range("E16").Value = '  +0.14%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '26.423.12'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.35%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '4.933'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.12%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.011'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.65%  '
$ws.Range("E20").Value = '  +0.89%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '190.95'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.23%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.241'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("E23").Value = '  +0.83%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '149.83'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1314'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +5.57%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.857'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.12%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.68'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.54%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06317'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.66%  '
$ws.Range("E29").Value = '  +4.69%  '
$ws.Range("E30").Value = '  -0.05%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.588'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.52%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.577'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.687'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.40%  '
$ws.Range("E34").Value = '  +0.76%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6161'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.88%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.414'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.711'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -0.39%  '
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01626'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.31%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.108.05'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8821'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.57%  '
$ws.Range("E42").Value = '  +0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '101.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.837.92'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.13%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00000000110'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +1.12%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '57.33'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.33%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.204'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.012'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.44%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05270'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.18%  '
$ws.Range("E50").Value = '  +0.37%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.044'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.23%  '
